# Fruta / hortaliza, semanal
# Insert a new weekly data row at row 86 (pushing the existing rows 86..161
# down to 87..162) and populate it with the new observation.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows(86).Insert()

$ws.Cells.Item(86, 1).Value = 8
$ws.Cells.Item(86, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(86, 3).Value = "Coquimbo"
$ws.Cells.Item(86, 4).Value = 44789
$ws.Cells.Item(86, 5).Value = 4
$ws.Cells.Item(86, 6).Value = 100112040
$ws.Cells.Item(86, 7).Value = "Cilantro"
$ws.Cells.Item(86, 8).Value = "Sin especificar"
$ws.Cells.Item(86, 9).Value = "Primera"
$ws.Cells.Item(86, 10).Value = 3200
$ws.Cells.Item(86, 11).Value = 2000
$ws.Cells.Item(86, 12).Value = 2500
$ws.Cells.Item(86, 13).Value = 2250
$ws.Cells.Item(86, 14).Value = "`$/atado 1 a 1,5 kilos"
$ws.Cells.Item(86, 15).Value = "Provincia del Elquí"
$ws.Cells.Item(86, 16).Value = 1500
$ws.Cells.Item(86, 17).Value = 1.5
$ws.Cells.Item(86, 18).Value = "Hortaliza"
